$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 17031.584
$ws.Range("I62").Value = 8118
$ws.Range("J62").Value = 61599.5
$ws.Range("K62").Value = 8118
$ws.Range("L62").Value = 61599.5
$ws.Range("M62").Value = -7494
$ws.Range("N62").Value = -62847.5

$ws.Range("H65").Value = 17031.584
$ws.Range("I65").Value = 8118
$ws.Range("J65").Value = 61599.5
$ws.Range("K65").Value = 40590
$ws.Range("L65").Value = 307997.5
$ws.Range("M65").Value = -37470
$ws.Range("N65").Value = -314237.5

$ws.Range("H69").Value = 250003150
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 250003150
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 750009450
$ws.Range("N69").Value = -750011198
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 250003150
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 250003150
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 2250028350
$ws.Range("N72").Value = -2250037086
$ws.Range("M72").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H132").Value = 2921.5
$ws.Range("I132").Value = 2711.0293
$ws.Range("J132").Value = 6499.5
$ws.Range("K132").Value = 8133.0879
$ws.Range("L132").Value = 19498.5
$ws.Range("M132").Value = -5603.0879
$ws.Range("N132").Value = -24558.5

$ws.Range("H140").Value = 71685.53
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 71685.53
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 71685.53
$ws.Range("N140").Value = -82045.53

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2385.4
$ws.Range("I45").Value = 2112.5
$ws.Range("J45").Value = 2697.2856
$ws.Range("K45").Value = 2112.5
$ws.Range("L45").Value = 2697.2856
$ws.Range("M45").Value = -1735.5
$ws.Range("N45").Value = -3451.2856

$ws.Range("H80").Value = 19499
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 19499
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 19499
$ws.Range("N80").Value = -21495

$ws.Range("H83").Value = 19499
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 19499
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 58497
$ws.Range("N83").Value = -68481

$ws.Range("H104").Value = 67573.336
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 67573.336
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 67573.336
$ws.Range("N104").Value = -74561.336

$ws.Range("H122").Value = 50500000
$ws.Range("I122").Value = 1000000
$ws.Range("J122").Value = 100000000
$ws.Range("K122").Value = 3000000
$ws.Range("L122").Value = 300000000
$ws.Range("M122").Value = -2997550
$ws.Range("N122").Value = -300004900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5400.7534
$ws.Range("I31").Value = 1764.9584
$ws.Range("J31").Value = 7181.551
$ws.Range("K31").Value = 1764.9584
$ws.Range("L31").Value = 7181.551
$ws.Range("M31").Value = -1469.9584
$ws.Range("N31").Value = -7771.551

$ws.Range("H34").Value = 5400.7534
$ws.Range("I34").Value = 1764.9584
$ws.Range("J34").Value = 7181.551
$ws.Range("K34").Value = 1764.9584
$ws.Range("L34").Value = 7181.551
$ws.Range("M34").Value = -1562.9584
$ws.Range("N34").Value = -7585.551

$ws.Range("H58").Value = 1669.2307
$ws.Range("I58").Value = 1300
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 1300
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -1097
$ws.Range("N58").Value = -2906

$ws.Range("H74").Value = 20191.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20191.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 20191.4
$ws.Range("N74").Value = -21939.4

$ws.Range("H77").Value = 20191.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20191.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 60574.2
$ws.Range("N77").Value = -69310.20000000001

$ws.Range("H122").Value = 1972.5555
$ws.Range("I122").Value = 1688.25
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 5064.75
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -2614.75
$ws.Range("N122").Value = -11500

$ws.Range("H134").Value = 16174.857
$ws.Range("I134").Value = 21644.8
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 64934.39999999999
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -62399.39999999999

$ws.Range("H136").Value = 1669.2307
$ws.Range("I136").Value = 1300
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 3900
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -1350
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1546.3334
$ws.Range("I55").Value = 1133.3334
$ws.Range("J55").Value = 1649.5834
$ws.Range("K55").Value = 3400.0002
$ws.Range("L55").Value = 4948.7502
$ws.Range("M55").Value = -3223.0002
$ws.Range("N55").Value = -5302.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 54665.832
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 54665.832
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 54665.832
$ws.Range("N82").Value = -55431.832

$ws.Range("H85").Value = 54665.832
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 54665.832
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 54665.832
$ws.Range("N85").Value = -57317.832

$ws.Range("H122").Value = 6116.7393
$ws.Range("I122").Value = 6116.7393
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18350.2179
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15900.2179
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2377.3
$ws.Range("I82").Value = 2967
$ws.Range("J82").Value = 2124.5715
$ws.Range("K82").Value = 2967
$ws.Range("L82").Value = 2124.5715
$ws.Range("M82").Value = -2606
$ws.Range("N82").Value = -2846.5715

$ws.Range("H85").Value = 2377.3
$ws.Range("I85").Value = 2967
$ws.Range("J85").Value = 2124.5715
$ws.Range("K85").Value = 2967
$ws.Range("L85").Value = 2124.5715
$ws.Range("M85").Value = -1719
$ws.Range("N85").Value = -4620.5715

$ws.Range("H94").Value = 76443.336
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 76443.336
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 76443.336
$ws.Range("N94").Value = -77795.336

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 2179.1333
$ws.Range("I122").Value = 2179.1333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6537.3999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4087.3999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2023.4667
$ws.Range("I122").Value = 2023.4667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6070.4001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3620.4001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3335645.2
$ws.Range("I132").Value = 1875.5927
$ws.Range("J132").Value = 7249201
$ws.Range("K132").Value = 5626.7781
$ws.Range("L132").Value = 21747603
$ws.Range("M132").Value = -3096.7781
$ws.Range("N132").Value = -21752663
